$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.42151153087616
$ws.Range("B1").Value = 1.657349944114685
$ws.Range("C1").Value = 5.115967750549316
$ws.Range("D1").Value = 2.84512996673584
$ws.Range("E1").Value = 0.871090829372406
